# "Generate Report for Handback" — mark the zh-cn/de-de handback as complete:
#  - Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" (shared by the zh-cn/de-de status cells)
#    and columns E/F widen to fit the new text.
#  - zh-cn / de-de sheets: fill in "Latest Target File" (I) + "Latest Handback
#    File" (J) with the handback markdown/link and, for de-de, refresh
#    "Latest Handback DateTime" (K) to the new timestamp.

$wb = $excel.ActiveWorkbook

$mdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2038e9f26585c2f1535a68cf53cd1d2e7132c5f/e2e/f2d62a4c-e913-49d1-af2b-2aba803b41bf.md"
$mdDisplay = "f2d62a4c-e913-49d1-af2b-2aba803b41bf.md"

# ---- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn sheet ---------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2038e9f26585c2f1535a68cf53cd1d2e7132c5f/e2e/f2d62a4c-e913-49d1-af2b-2aba803b41bf.md", "", "", "f2d62a4c-e913-49d1-af2b-2aba803b41bf.md") | Out-Null
$zhcn.Range("J2").Value = "f2d62a4c-e913-49d1-af2b-2aba803b41bf.06ad2a0f9e6932d06e5ae424372b6fa2f1326765.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdTarget, "", "", $mdDisplay) | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2038e9f26585c2f1535a68cf53cd1d2e7132c5f/e2e/ffffe0181847-0564-4a35-b1c6-304a42343a77.md", "", "", "ffffe0181847-0564-4a35-b1c6-304a42343a77.md") | Out-Null
$zhcn.Range("J3").Value = "f2d62a4c-e913-49d1-af2b-2aba803b41bf.06ad2a0f9e6932d06e5ae424372b6fa2f1326765.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdTarget, "", "", $mdDisplay) | Out-Null

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2038e9f26585c2f1535a68cf53cd1d2e7132c5f/e2e/f2d62a4c-e913-49d1-af2b-2aba803b41bf.md", "", "", "f2d62a4c-e913-49d1-af2b-2aba803b41bf.md") | Out-Null
$dede.Range("J2").Value = "f2d62a4c-e913-49d1-af2b-2aba803b41bf.06ad2a0f9e6932d06e5ae424372b6fa2f1326765.de-de.xlf"
$dede.Range("K2").Value = "2016-08-17 13:00:55"
$dede.Hyperlinks.Add($dede.Range("I2"), $mdTarget, "", "", $mdDisplay) | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2038e9f26585c2f1535a68cf53cd1d2e7132c5f/e2e/ffffe0181847-0564-4a35-b1c6-304a42343a77.md", "", "", "ffffe0181847-0564-4a35-b1c6-304a42343a77.md") | Out-Null
$dede.Range("J3").Value = "f2d62a4c-e913-49d1-af2b-2aba803b41bf.06ad2a0f9e6932d06e5ae424372b6fa2f1326765.de-de.xlf"
$dede.Range("K3").Value = "2016-08-17 13:00:55"
$dede.Hyperlinks.Add($dede.Range("I3"), $mdTarget, "", "", $mdDisplay) | Out-Null

# ---- shared string for the zh-cn "Latest Handback DateTime" -------------
$zhcn.Range("K2").Value = "2016-08-17 13:00:48"
$zhcn.Range("K3").Value = "2016-08-17 13:00:48"
